$wb = $excel.ActiveWorkbook

# Existing sheets before edit: "ODI Batting" (sheet1), "ODI Bowling" (sheet2)
$battingWs = $wb.Worksheets.Item("ODI Batting")

# --- Add new "Player Info" sheet, inserted before "ODI Batting" ---
$infoWs = $wb.Worksheets.Add($battingWs)
$infoWs.Name = "Player Info"

# Re-fetch sheets by name now that the workbook structure changed, to avoid stale references
$infoWs = $wb.Worksheets.Item("Player Info")
$battingWs = $wb.Worksheets.Item("ODI Batting")
$bowlingWs = $wb.Worksheets.Item("ODI Bowling")

# Header row for "Player Info"
$infoWs.Range("A1").Value = "ID"
$infoWs.Range("B1").Value = "NAME"
$infoWs.Range("C1").Value = "BATTING_HAND"
$infoWs.Range("D1").Value = "BOWL_STYLE"

# Match the header formatting used on the other sheets (bold, centered/top aligned, thin box border)
$infoWs.Range("A1:D1").Font.Bold = $true
$infoWs.Range("A1:D1").HorizontalAlignment = -4108
$infoWs.Range("A1:D1").VerticalAlignment = -4160
$infoWs.Range("A1:D1").Borders.LineStyle = 1

# Data row for "Player Info" -- keep ID as text (matches source data convention of text cells)
$infoWs.Range("A2").NumberFormat = "@"
$infoWs.Range("A2").Value = "4590"
$infoWs.Range("A2").Style = "Normal"
$infoWs.Range("B2").Value = "Naveen-ul-Haq Murid"
$infoWs.Range("C2").Value = "Right Handed"
$infoWs.Range("D2").Value = "Right Arm Medium Fast"

# --- Update "ODI Batting" sheet: MATCH_CARD_LINK -> MATCH_CODE, full URL -> bare match code ---
$battingWs.Range("D1").Value = "MATCH_CODE"

$battingCodes = $battingWs.Range("D2:D8")
$battingCodes.NumberFormat = "@"
$battingWs.Range("D2").Value = "3936"
$battingWs.Range("D3").Value = "3938"
$battingWs.Range("D4").Value = "4377"
$battingWs.Range("D5").Value = "4378"
$battingWs.Range("D6").Value = "4444"
$battingWs.Range("D7").Value = "4446"
$battingWs.Range("D8").Value = "4448"
$battingCodes.Style = "Normal"

# --- Update "ODI Bowling" sheet: MATCH_CARD_LINK -> MATCH_CODE, full URL -> bare match code ---
$bowlingWs.Range("B1").Value = "MATCH_CODE"

$bowlingCodes = $bowlingWs.Range("B2:B8")
$bowlingCodes.NumberFormat = "@"
$bowlingWs.Range("B2").Value = "3936"
$bowlingWs.Range("B3").Value = "3938"
$bowlingWs.Range("B4").Value = "4377"
$bowlingWs.Range("B5").Value = "4378"
$bowlingWs.Range("B6").Value = "4444"
$bowlingWs.Range("B7").Value = "4446"
$bowlingWs.Range("B8").Value = "4448"
$bowlingCodes.Style = "Normal"
